$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 26468
$ws.Range("E2").Value = 512903223779
$ws.Range("F2").Value = 7126326631
$ws.Range("G2").Value = 0.36415
$ws.Range("D3").Value = 1808.34
$ws.Range("E3").Value = 217450894711
$ws.Range("F3").Value = 7524970555
$ws.Range("G3").Value = 0.41
$ws.Range("E4").Value = 83043328337
$ws.Range("F4").Value = 11334795825
$ws.Range("G4").Value = 0.06443
$ws.Range("D5").Value = 304.6
$ws.Range("E5").Value = 48090100899
$ws.Range("F5").Value = 551217400
$ws.Range("G5").Value = -0.49312
$ws.Range("D6").Value = 0.999892
$ws.Range("E6").Value = 29054489055
$ws.Range("F6").Value = 3737060407
$ws.Range("G6").Value = -0.11316
$ws.Range("D7").Value = 0.453731
$ws.Range("E7").Value = 23537286259
$ws.Range("F7").Value = 824227008
$ws.Range("G7").Value = -0.02674
$ws.Range("D8").Value = 0.358471
$ws.Range("E8").Value = 12560534006
$ws.Range("F8").Value = 173699781
$ws.Range("G8").Value = -1.70864
$ws.Range("D9").Value = 1806.75
$ws.Range("E9").Value = 12083645961
$ws.Range("F9").Value = 7192764
$ws.Range("G9").Value = 0.26433
$ws.Range("D10").Value = 0.07094
$ws.Range("E10").Value = 9892292213
$ws.Range("F10").Value = 277282380
$ws.Range("G10").Value = -0.2264
$ws.Range("D11").Value = 0.888933
$ws.Range("E11").Value = 8246329628
$ws.Range("F11").Value = 293156893
$ws.Range("G11").Value = 1.58333
$ws.Range("D12").Value = 19.25
$ws.Range("E12").Value = 7623017389
$ws.Range("F12").Value = 235275661
$ws.Range("G12").Value = -0.12426
$ws.Range("D13").Value = 0.076961
$ws.Range("E13").Value = 6951799441
$ws.Range("F13").Value = 364631544
$ws.Range("G13").Value = -0.23385
$ws.Range("D14").Value = 5.24
$ws.Range("E14").Value = 6482169999
$ws.Range("F14").Value = 108562153
$ws.Range("G14").Value = -0.40531
$ws.Range("D15").Value = 85.54000000000001
$ws.Range("E15").Value = 6247949254
$ws.Range("F15").Value = 632241522
$ws.Range("G15").Value = -0.44344
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 5302793567
$ws.Range("F16").Value = 1509861258
$ws.Range("G16").Value = -0.06693
$ws.Range("D17").Value = 0.00000851
$ws.Range("E17").Value = 5010898637
$ws.Range("F17").Value = 132071460
$ws.Range("G17").Value = -0.56887
$ws.Range("D18").Value = 14.09
$ws.Range("E18").Value = 4711229912
$ws.Range("F18").Value = 148802188
$ws.Range("G18").Value = -0.48773
$ws.Range("D19").Value = 0.99991
$ws.Range("E19").Value = 4621833351
$ws.Range("F19").Value = 94486427
$ws.Range("G19").Value = -0.07106999999999999
$ws.Range("D20").Value = 26537
$ws.Range("E20").Value = 4138886776
$ws.Range("F20").Value = 81966498
$ws.Range("G20").Value = 0.35829
$ws.Range("D21").Value = 4.94
$ws.Range("E21").Value = 3718359137
$ws.Range("F21").Value = 35995431
$ws.Range("G21").Value = -0.72487
$ws.Range("E22").Value = 3308003035
$ws.Range("F22").Value = 211594
$ws.Range("G22").Value = 0.23301
$ws.Range("E23").Value = 3246793498
$ws.Range("F23").Value = 121205122
$ws.Range("G23").Value = -0.93057
$ws.Range("E24").Value = 3062642730
$ws.Range("F24").Value = 73716403
$ws.Range("G24").Value = 0.57355
$ws.Range("E25").Value = 2823586431
$ws.Range("F25").Value = 8108704
$ws.Range("G25").Value = -2.64663
$ws.Range("D26").Value = 46.15
$ws.Range("E26").Value = 2772085020
$ws.Range("F26").Value = 16856007
$ws.Range("G26").Value = 2.98342
$ws.Range("D27").Value = 151.5
$ws.Range("E27").Value = 2749388901
$ws.Range("F27").Value = 75473165
$ws.Range("G27").Value = 0.53229
$ws.Range("E28").Value = 2499127273
$ws.Range("F28").Value = 70673575
$ws.Range("G28").Value = -0.59949
$ws.Range("D29").Value = 0.086981
$ws.Range("E29").Value = 2331041613
$ws.Range("F29").Value = 32596627
$ws.Range("G29").Value = 0.5602200000000001
$ws.Range("D30").Value = 111.92
$ws.Range("E30").Value = 2172073269
$ws.Range("F30").Value = 56907987
$ws.Range("G30").Value = -0.37054
$ws.Range("E31").Value = 2095066658
$ws.Range("F31").Value = 23939677
$ws.Range("G31").Value = -0.67116
$ws.Range("D32").Value = 0.999327
$ws.Range("E32").Value = 2040783273
$ws.Range("F32").Value = 201451772
$ws.Range("G32").Value = -0.15587
$ws.Range("D33").Value = 4.41
$ws.Range("E33").Value = 1882967344
$ws.Range("F33").Value = 86976904
$ws.Range("G33").Value = -1.5852
$ws.Range("E34").Value = 1771653517
$ws.Range("F34").Value = 58775625
$ws.Range("G34").Value = -1.83161
$ws.Range("D35").Value = 0.050867
$ws.Range("E35").Value = 1600267163
$ws.Range("F35").Value = 16283866
$ws.Range("G35").Value = -0.47853
$ws.Range("E36").Value = 1585039158
$ws.Range("F36").Value = 40303449
$ws.Range("G36").Value = -1.48674
$ws.Range("D37").Value = 0.059875
$ws.Range("E37").Value = 1512806160
$ws.Range("F37").Value = 5291641
$ws.Range("G37").Value = -0.33828
$ws.Range("D38").Value = 98.54000000000001
$ws.Range("E38").Value = 1432057912
$ws.Range("F38").Value = 18021819
$ws.Range("G38").Value = -2.01079
$ws.Range("E39").Value = 1414516098
$ws.Range("F39").Value = 58764707
$ws.Range("G39").Value = -1.67306
$ws.Range("E40").Value = 1410788169
$ws.Range("F40").Value = 167307036
$ws.Range("G40").Value = -0.72015
$ws.Range("D41").Value = 0.01930521
$ws.Range("E41").Value = 1403540620
$ws.Range("F41").Value = 38554262
$ws.Range("G41").Value = -0.51497
$ws.Range("E42").Value = 1191680697
$ws.Range("F42").Value = 63347754
$ws.Range("G42").Value = -3.55297
$ws.Range("D43").Value = 0.150737
$ws.Range("E43").Value = 1092161362
$ws.Range("F43").Value = 65689689
$ws.Range("G43").Value = -3.27491
$ws.Range("B44").Value = "GGTKN"
$ws.Range("C44").Value = "GGTKN"
$ws.Range("D44").Value = 0.09175899999999999
$ws.Range("E44").Value = 1054238259
$ws.Range("F44").Value = 68561
$ws.Range("G44").Value = 1.4648
$ws.Range("B45").Value = "GRT"
$ws.Range("C45").Value = "The Graph"
$ws.Range("D45").Value = 0.116756
$ws.Range("E45").Value = 1049367196
$ws.Range("F45").Value = 31441752
$ws.Range("G45").Value = 0.2163
$ws.Range("D46").Value = 0.999765
$ws.Range("E46").Value = 1018788196
$ws.Range("F46").Value = 27999688
$ws.Range("G46").Value = -0.05818
$ws.Range("D47").Value = 0.9993649999999999
$ws.Range("E47").Value = 1003111312
$ws.Range("F47").Value = 15645282
$ws.Range("G47").Value = -0.09135
$ws.Range("D48").Value = 0.999676
$ws.Range("E48").Value = 999875549
$ws.Range("F48").Value = 5368823
$ws.Range("G48").Value = 0.00346
$ws.Range("E49").Value = 992159152
$ws.Range("F49").Value = 222893434
$ws.Range("G49").Value = 1.49577
$ws.Range("D50").Value = 0.5086000000000001
$ws.Range("E50").Value = 942814823
$ws.Range("F50").Value = 174324609
$ws.Range("G50").Value = 1.82811
$ws.Range("D51").Value = 0.843603
$ws.Range("E51").Value = 934024033
$ws.Range("F51").Value = 84618218
$ws.Range("G51").Value = 0.50784
